$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.829.56'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '3.381.06'
$ws.Range('E3').Value = '  +7.89%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'" + '261.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.67%  '
$ws.Range('D6').Value = "'" + '634.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('D7').Value = "'" + '1.38'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +23.44%  '
$ws.Range('D8').Value = "'" + '0.396'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').Value = "'" + '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = "'" + '0.883'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.76%  '
$ws.Range('D11').Value = '3.380.07'
$ws.Range('E11').Value = '  +8.00%  '
$ws.Range('D12').Value = "'" + '0.200'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').Value = '98.651.29'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = "'" + '36.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('D15').Value = "'" + '0.0000250'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = '4.005.99'
$ws.Range('E16').Value = '  +7.92%  '
$ws.Range('D17').Value = "'" + '5.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '3.375.38'
$ws.Range('E18').Value = '  +8.57%  '
$ws.Range('D19').Value = "'" + '3.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').Value = "'" + '15.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.84%  '
$ws.Range('D21').Value = "'" + '496.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.10%  '
$ws.Range('E22').Value = '  +8.70%  '
$ws.Range('D23').Value = "'" + '0.0000212'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.45%  '
$ws.Range('D24').Value = "'" + '9.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.99%  '
$ws.Range('D25').Value = "'" + '5.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').Value = "'" + '89.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('D27').Value = "'" + '12.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.98%  '
$ws.Range('D28').Value = '3.514.26'
$ws.Range('E28').Value = '  +6.81%  '
$ws.Range('D29').Value = "'" + '0.283'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.72%  '
$ws.Range('B30').Value = 'Cronos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D30').Value = "'" + '0.199'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +12.97%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value = "'" + '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').Value = "'" + '0.135'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.66%  '
$ws.Range('D33').Value = "'" + '9.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.55%  '
$ws.Range('E34').Value = '  +18.84%  '
$ws.Range('D35').Value = "'" + '27.91'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.80%  '
$ws.Range('D36').Value = "'" + '7.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'" + '0.151'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').Value = "'" + '1.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.81%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = "'" + '507.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.22%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = "'" + '0.471'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.00%  '
$ws.Range('D41').Value = "'" + '24.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = "'" + '1.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('B43').Value = 'MantraDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D43').Value = "'" + '3.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('D44').Value = "'" + '3.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.67%  '
$ws.Range('D45').Value = "'" + '0.788'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.98%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = "'" + '160.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('E49').Value = '  +6.54%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = "'" + '0.838'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.83%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = "'" + '46.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.55%  '
